# Update cryptocurrency market data (Daten aktualisiert am 2024-04-07)
# Several coins changed rank order (B/C ticker+name reshuffled) and
# Price / Market Cap / Volume / Change(24h) columns were refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 69306
$ws.Cells.Item(2, 5).Value = 1363932763074
$ws.Cells.Item(2, 6).Value = 19008976614
$ws.Cells.Item(2, 7).Value = 2.3426
$ws.Cells.Item(3, 4).Value = 3387.36
$ws.Cells.Item(3, 5).Value = 406710446312
$ws.Cells.Item(3, 6).Value = 8724892660
$ws.Cells.Item(3, 7).Value = 1.62731
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 106875023058
$ws.Cells.Item(4, 6).Value = 28937788824
$ws.Cells.Item(4, 7).Value = -0.05547
$ws.Cells.Item(5, 4).Value = 586.52
$ws.Cells.Item(5, 5).Value = 90225732038
$ws.Cells.Item(5, 6).Value = 696239279
$ws.Cells.Item(5, 7).Value = 0.93176
$ws.Cells.Item(6, 4).Value = 180.46
$ws.Cells.Item(6, 5).Value = 80447359878
$ws.Cells.Item(6, 6).Value = 2243774447
$ws.Cells.Item(6, 7).Value = 2.49275
$ws.Cells.Item(7, 2).Value = "USDC"
$ws.Cells.Item(7, 3).Value = "USDC"
$ws.Cells.Item(7, 4).Value = 0.999366
$ws.Cells.Item(7, 5).Value = 32910098254
$ws.Cells.Item(7, 6).Value = 3919859957
$ws.Cells.Item(7, 7).Value = -0.05639
$ws.Cells.Item(8, 4).Value = 0.5969370000000001
$ws.Cells.Item(8, 5).Value = 32859797063
$ws.Cells.Item(8, 6).Value = 806482083
$ws.Cells.Item(8, 7).Value = 1.48356
$ws.Cells.Item(9, 2).Value = "STETH"
$ws.Cells.Item(9, 3).Value = "Lido Staked Ether"
$ws.Cells.Item(9, 4).Value = 3381.62
$ws.Cells.Item(9, 5).Value = 32239574756
$ws.Cells.Item(9, 6).Value = 101616432
$ws.Cells.Item(9, 7).Value = 1.64889
$ws.Cells.Item(10, 4).Value = 0.197361
$ws.Cells.Item(10, 5).Value = 28422027126
$ws.Cells.Item(10, 6).Value = 2356790810
$ws.Cells.Item(10, 7).Value = 8.09686
$ws.Cells.Item(11, 4).Value = 0.593552
$ws.Cells.Item(11, 5).Value = 20938694026
$ws.Cells.Item(11, 6).Value = 257759748
$ws.Cells.Item(11, 7).Value = 2.28471
$ws.Cells.Item(12, 2).Value = "TON"
$ws.Cells.Item(12, 3).Value = "Toncoin"
$ws.Cells.Item(12, 4).Value = 5.41
$ws.Cells.Item(12, 5).Value = 18755291051
$ws.Cells.Item(12, 6).Value = 142795156
$ws.Cells.Item(12, 7).Value = -2.08733
$ws.Cells.Item(13, 2).Value = "AVAX"
$ws.Cells.Item(13, 3).Value = "Avalanche"
$ws.Cells.Item(13, 4).Value = 48.5
$ws.Cells.Item(13, 5).Value = 18332218378
$ws.Cells.Item(13, 6).Value = 492988753
$ws.Cells.Item(13, 7).Value = 3.38384
$ws.Cells.Item(14, 4).Value = 0.00002834
$ws.Cells.Item(14, 5).Value = 16702587269
$ws.Cells.Item(14, 6).Value = 603249154
$ws.Cells.Item(14, 7).Value = 4.26076
$ws.Cells.Item(15, 2).Value = "BCH"
$ws.Cells.Item(15, 3).Value = "Bitcoin Cash"
$ws.Cells.Item(15, 4).Value = 679.48
$ws.Cells.Item(15, 5).Value = 13382133896
$ws.Cells.Item(15, 6).Value = 773957282
$ws.Cells.Item(15, 7).Value = -1.89252
$ws.Cells.Item(16, 2).Value = "DOT"
$ws.Cells.Item(16, 3).Value = "Polkadot"
$ws.Cells.Item(16, 4).Value = 8.66
$ws.Cells.Item(16, 5).Value = 11684042795
$ws.Cells.Item(16, 6).Value = 169388464
$ws.Cells.Item(16, 7).Value = 3.06631
$ws.Cells.Item(17, 2).Value = "WBTC"
$ws.Cells.Item(17, 3).Value = "Wrapped Bitcoin"
$ws.Cells.Item(17, 4).Value = 69333
$ws.Cells.Item(17, 5).Value = 10779832953
$ws.Cells.Item(17, 6).Value = 141651025
$ws.Cells.Item(17, 7).Value = 2.34724
$ws.Cells.Item(18, 2).Value = "TRX"
$ws.Cells.Item(18, 3).Value = "TRON"
$ws.Cells.Item(18, 4).Value = 0.120311
$ws.Cells.Item(18, 5).Value = 10550006145
$ws.Cells.Item(18, 6).Value = 225653074
$ws.Cells.Item(18, 7).Value = 1.46164
$ws.Cells.Item(19, 2).Value = "LINK"
$ws.Cells.Item(19, 3).Value = "Chainlink"
$ws.Cells.Item(19, 4).Value = 17.73
$ws.Cells.Item(19, 5).Value = 10411138585
$ws.Cells.Item(19, 6).Value = 230512399
$ws.Cells.Item(19, 7).Value = 0.74143
$ws.Cells.Item(20, 4).Value = 11.32
$ws.Cells.Item(20, 5).Value = 8520482254
$ws.Cells.Item(20, 6).Value = 157968113
$ws.Cells.Item(20, 7).Value = 2.34232
$ws.Cells.Item(21, 4).Value = 0.906156
$ws.Cells.Item(21, 5).Value = 8413234951
$ws.Cells.Item(21, 6).Value = 221104999
$ws.Cells.Item(21, 7).Value = 1.51905
$ws.Cells.Item(22, 4).Value = 17.2
$ws.Cells.Item(22, 5).Value = 7955984665
$ws.Cells.Item(22, 6).Value = 121881233
$ws.Cells.Item(22, 7).Value = 1.8232
$ws.Cells.Item(23, 4).Value = 103.18
$ws.Cells.Item(23, 5).Value = 7677385209
$ws.Cells.Item(23, 6).Value = 601614109
$ws.Cells.Item(23, 7).Value = 1.83649
$ws.Cells.Item(24, 4).Value = 6.95
$ws.Cells.Item(24, 5).Value = 7360388208
$ws.Cells.Item(24, 6).Value = 367999243
$ws.Cells.Item(24, 7).Value = -1.5515
$ws.Cells.Item(25, 4).Value = 13.76
$ws.Cells.Item(25, 5).Value = 5481145643
$ws.Cells.Item(25, 6).Value = 180712599
$ws.Cells.Item(25, 7).Value = 5.86007
$ws.Cells.Item(26, 4).Value = 5.78
$ws.Cells.Item(26, 5).Value = 5350892072
$ws.Cells.Item(26, 6).Value = 1273179
$ws.Cells.Item(26, 7).Value = 0.23891
$ws.Cells.Item(27, 2).Value = "DAI"
$ws.Cells.Item(27, 3).Value = "Dai"
$ws.Cells.Item(27, 4).Value = 0.999645
$ws.Cells.Item(27, 5).Value = 5106876172
$ws.Cells.Item(27, 6).Value = 646973416
$ws.Cells.Item(27, 7).Value = -0.13723
$ws.Cells.Item(28, 2).Value = "ETC"
$ws.Cells.Item(28, 3).Value = "Ethereum Classic"
$ws.Cells.Item(28, 4).Value = 33.86
$ws.Cells.Item(28, 5).Value = 4960004530
$ws.Cells.Item(28, 6).Value = 206460816
$ws.Cells.Item(28, 7).Value = 2.04751
$ws.Cells.Item(29, 2).Value = "STX"
$ws.Cells.Item(29, 3).Value = "Stacks"
$ws.Cells.Item(29, 4).Value = 3.27
$ws.Cells.Item(29, 5).Value = 4752287769
$ws.Cells.Item(29, 6).Value = 57154795
$ws.Cells.Item(29, 7).Value = 3.63692
$ws.Cells.Item(30, 2).Value = "FIL"
$ws.Cells.Item(30, 3).Value = "Filecoin"
$ws.Cells.Item(30, 4).Value = 8.789999999999999
$ws.Cells.Item(30, 5).Value = 4687935910
$ws.Cells.Item(30, 6).Value = 176177202
$ws.Cells.Item(30, 7).Value = 2.88181
$ws.Cells.Item(31, 2).Value = "MNT"
$ws.Cells.Item(31, 3).Value = "Mantle"
$ws.Cells.Item(31, 4).Value = 1.4
$ws.Cells.Item(31, 5).Value = 4457565465
$ws.Cells.Item(31, 6).Value = 42694505
$ws.Cells.Item(31, 7).Value = 5.31252
$ws.Cells.Item(32, 2).Value = "ATOM"
$ws.Cells.Item(32, 3).Value = "Cosmos Hub"
$ws.Cells.Item(32, 4).Value = 11.18
$ws.Cells.Item(32, 5).Value = 4367259087
$ws.Cells.Item(32, 6).Value = 170447907
$ws.Cells.Item(32, 7).Value = 1.56546
$ws.Cells.Item(33, 4).Value = 1.5
$ws.Cells.Item(33, 5).Value = 3969519339
$ws.Cells.Item(33, 6).Value = 168585528
$ws.Cells.Item(33, 7).Value = 0.7558
$ws.Cells.Item(34, 2).Value = "CRO"
$ws.Cells.Item(34, 3).Value = "Cronos"
$ws.Cells.Item(34, 4).Value = 0.146419
$ws.Cells.Item(34, 5).Value = 3899929131
$ws.Cells.Item(34, 6).Value = 10582675
$ws.Cells.Item(34, 7).Value = 2.32239
$ws.Cells.Item(35, 2).Value = "IMX"
$ws.Cells.Item(35, 3).Value = "Immutable"
$ws.Cells.Item(35, 4).Value = 2.73
$ws.Cells.Item(35, 5).Value = 3881452208
$ws.Cells.Item(35, 6).Value = 24725941
$ws.Cells.Item(35, 7).Value = 1.82722
$ws.Cells.Item(36, 2).Value = "XLM"
$ws.Cells.Item(36, 3).Value = "Stellar"
$ws.Cells.Item(36, 4).Value = 0.130072
$ws.Cells.Item(36, 5).Value = 3752192200
$ws.Cells.Item(36, 6).Value = 48970072
$ws.Cells.Item(36, 7).Value = 1.29402
$ws.Cells.Item(37, 2).Value = "RNDR"
$ws.Cells.Item(37, 3).Value = "Render"
$ws.Cells.Item(37, 4).Value = 9.710000000000001
$ws.Cells.Item(37, 5).Value = 3707584587
$ws.Cells.Item(37, 6).Value = 87759660
$ws.Cells.Item(37, 7).Value = 3.11429
$ws.Cells.Item(38, 2).Value = "TAO"
$ws.Cells.Item(38, 3).Value = "Bittensor"
$ws.Cells.Item(38, 4).Value = 556.42
$ws.Cells.Item(38, 5).Value = 3648411660
$ws.Cells.Item(38, 6).Value = 13830418
$ws.Cells.Item(38, 7).Value = -2.52961
$ws.Cells.Item(39, 2).Value = "WIF"
$ws.Cells.Item(39, 3).Value = "dogwifhat"
$ws.Cells.Item(39, 4).Value = 3.62
$ws.Cells.Item(39, 5).Value = 3618850501
$ws.Cells.Item(39, 6).Value = 392057021
$ws.Cells.Item(39, 7).Value = 10.73044
$ws.Cells.Item(40, 2).Value = "HBAR"
$ws.Cells.Item(40, 3).Value = "Hedera"
$ws.Cells.Item(40, 4).Value = 0.106466
$ws.Cells.Item(40, 5).Value = 3592627774
$ws.Cells.Item(40, 6).Value = 25870450
$ws.Cells.Item(40, 7).Value = 1.27874
$ws.Cells.Item(41, 2).Value = "FDUSD"
$ws.Cells.Item(41, 3).Value = "First Digital USD"
$ws.Cells.Item(41, 4).Value = 1.002
$ws.Cells.Item(41, 5).Value = 3530854261
$ws.Cells.Item(41, 6).Value = 4366880686
$ws.Cells.Item(41, 7).Value = 0.02028
$ws.Cells.Item(42, 2).Value = "OKB"
$ws.Cells.Item(42, 3).Value = "OKB"
$ws.Cells.Item(42, 4).Value = 58.55
$ws.Cells.Item(42, 5).Value = 3512831253
$ws.Cells.Item(42, 6).Value = 8389308
$ws.Cells.Item(42, 7).Value = 2.09003
$ws.Cells.Item(43, 2).Value = "MKR"
$ws.Cells.Item(43, 3).Value = "Maker"
$ws.Cells.Item(43, 4).Value = 3667.29
$ws.Cells.Item(43, 5).Value = 3390829114
$ws.Cells.Item(43, 6).Value = 91315242
$ws.Cells.Item(43, 7).Value = -1.07582
$ws.Cells.Item(44, 2).Value = "KAS"
$ws.Cells.Item(44, 3).Value = "Kaspa"
$ws.Cells.Item(44, 4).Value = 0.139934
$ws.Cells.Item(44, 5).Value = 3230871716
$ws.Cells.Item(44, 6).Value = 51320302
$ws.Cells.Item(44, 7).Value = 5.4082
$ws.Cells.Item(45, 2).Value = "GRT"
$ws.Cells.Item(45, 3).Value = "The Graph"
$ws.Cells.Item(45, 4).Value = 0.33975
$ws.Cells.Item(45, 5).Value = 3215485222
$ws.Cells.Item(45, 6).Value = 74158652
$ws.Cells.Item(45, 7).Value = 1.61073
$ws.Cells.Item(46, 2).Value = "INJ"
$ws.Cells.Item(46, 3).Value = "Injective"
$ws.Cells.Item(46, 4).Value = 35.56
$ws.Cells.Item(46, 5).Value = 3202455945
$ws.Cells.Item(46, 6).Value = 120757229
$ws.Cells.Item(46, 7).Value = 0.7668
$ws.Cells.Item(47, 2).Value = "VET"
$ws.Cells.Item(47, 3).Value = "VeChain"
$ws.Cells.Item(47, 4).Value = 0.04243749
$ws.Cells.Item(47, 5).Value = 3085243824
$ws.Cells.Item(47, 6).Value = 53112420
$ws.Cells.Item(47, 7).Value = 4.27306
$ws.Cells.Item(48, 2).Value = "OP"
$ws.Cells.Item(48, 3).Value = "Optimism"
$ws.Cells.Item(48, 4).Value = 3.06
$ws.Cells.Item(48, 5).Value = 3078515086
$ws.Cells.Item(48, 6).Value = 156220264
$ws.Cells.Item(48, 7).Value = 2.9313
$ws.Cells.Item(49, 2).Value = "PEPE"
$ws.Cells.Item(49, 3).Value = "Pepe"
$ws.Cells.Item(49, 4).Value = 0.00000719
$ws.Cells.Item(49, 5).Value = 3022980250
$ws.Cells.Item(49, 6).Value = 364608070
$ws.Cells.Item(49, 7).Value = 7.50048
$ws.Cells.Item(50, 2).Value = "FET"
$ws.Cells.Item(50, 3).Value = "Fetch.ai"
$ws.Cells.Item(50, 4).Value = 2.69
$ws.Cells.Item(50, 5).Value = 2809851322
$ws.Cells.Item(50, 6).Value = 175702584
$ws.Cells.Item(50, 7).Value = 3.24809
$ws.Cells.Item(51, 4).Value = 2.69
$ws.Cells.Item(51, 5).Value = 2691882686
$ws.Cells.Item(51, 6).Value = 28470429
$ws.Cells.Item(51, 7).Value = 1.54964
